# The document contains, in four places, an <id>...</id> tag split across
# three separate runs (e.g. "<id>" / "p076v_1" / "</id>"), left over from an
# earlier paste. Collapse each trio back into a single run holding the full
# "<id>p076v_N</id>" text (Word's Find/Replace naturally merges the runs it
# rewrites, taking on the formatting of the first one).
$d = $word.ActiveDocument

$ids = @("p076v_1", "p076v_2", "p076v_3", "p076v_4")
foreach ($id in $ids) {
    $tag = "<id>" + $id + "</id>"
    [void]$d.Content.Find.Execute($tag, $false, $false, $false, $false, $false, $true, 1, $false, $tag, 2)
}
